# Edit script: insert a new "side" column (camera placement direction)
# before the existing tool_anvil/tool_site/streambed columns, shifting
# them right by one column, then populate the new column with the
# appropriate categorical direction value for each camera row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F; old F (tool_anvil), G (tool_site), H
# (streambed) and I (notes) shift one column to the right automatically.
$ws.Columns("F:F").Insert()

# Header for the newly inserted column
$ws.Range("F1").Value = "side"

# side = west
$westRows = @(10, 11, 12, 22, 23, 24, 25, 26, 27, 33, 34, 35, 36, 40, 41, 42, 43, 45, 46, 49, 50, 51, 52, 57, 58, 59, 60, 61, 62, 63, 64, 65, 66, 67, 68, 69, 70, 71, 72, 73, 74, 75, 76, 77, 78, 79, 80, 81, 82, 83, 84, 116, 117, 118, 119, 120, 121, 122, 123, 124, 125, 126, 127, 128, 129, 130, 131, 132, 133, 134, 135, 136, 137, 138, 139, 140, 141, 142, 143)
foreach ($r in $westRows) {
    $ws.Cells.Item($r, 6).Value = "west"
}

# side = north
$northRows = @(2, 8, 15, 29, 30, 38, 39, 44, 47, 48, 144)
foreach ($r in $northRows) {
    $ws.Cells.Item($r, 6).Value = "north"
}

# side = east
$eastRows = @(13, 14, 16, 17, 28, 90, 91, 92, 93, 94, 95, 96, 97, 98, 99, 100, 101, 102, 103, 104, 105, 106, 107, 108, 109, 110, 111, 112, 113, 114, 115)
foreach ($r in $eastRows) {
    $ws.Cells.Item($r, 6).Value = "east"
}

# side = southeast
$southeastRows = @(3, 4, 37)
foreach ($r in $southeastRows) {
    $ws.Cells.Item($r, 6).Value = "southeast"
}

# side = south
$southRows = @(5, 6, 7, 9, 18, 19, 32)
foreach ($r in $southRows) {
    $ws.Cells.Item($r, 6).Value = "south"
}

# side = northwest
$northwestRows = @(20, 21, 53, 54, 55, 56, 85, 86, 87, 88, 89)
foreach ($r in $northwestRows) {
    $ws.Cells.Item($r, 6).Value = "northwest"
}

# Restore the saved view state: active selection
$ws.Range("F31").Select()
